$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "⭐Strong Buy" -> "Stock ⭐Strong Buy"
#    Insert a new leading column ("类别" / "股票") ahead of the ticker column.
# ---------------------------------------------------------------------------
$wsStockStrong = $wb.Worksheets.Item(1)
$wsStockStrong.Name = "Stock ⭐Strong Buy"

$wsStockStrong.Columns.Item(1).Insert()
$wsStockStrong.Range("B1").Copy() | Out-Null
$wsStockStrong.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$wsStockStrong.Range("A1").Value = "类别"
$wsStockStrong.Range("B1").Value = "代码"
for ($r = 2; $r -le 3; $r++) {
    $wsStockStrong.Range("A$r").Value = "股票"
}

# ---------------------------------------------------------------------------
# 2) "✅Buy" -> "Stock ✅Buy"
#    Same leading-column restructure.
# ---------------------------------------------------------------------------
$wsStockBuy = $wb.Worksheets.Item(2)
$wsStockBuy.Name = "Stock ✅Buy"

$wsStockBuy.Columns.Item(1).Insert()
$wsStockBuy.Range("B1").Copy() | Out-Null
$wsStockBuy.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$wsStockBuy.Range("A1").Value = "类别"
$wsStockBuy.Range("B1").Value = "代码"
for ($r = 2; $r -le 6; $r++) {
    $wsStockBuy.Range("A$r").Value = "股票"
}

# ---------------------------------------------------------------------------
# 3) New sheet "ETF ⭐Strong Buy" (headers only, right after "Stock ✅Buy")
# ---------------------------------------------------------------------------
$wsEtfStrong = $wb.Worksheets.Add($null, $wsStockBuy)
$wsEtfStrong.Name = "ETF ⭐Strong Buy"
$wsStockStrong.Range("A1").Copy() | Out-Null
$wsEtfStrong.Range("A1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$wsEtfStrong.Range("A1").Value = "类别"
$wsEtfStrong.Range("B1").Value = "代码"
$wsEtfStrong.Range("C1").Value = "收盘价"
$wsEtfStrong.Range("D1").Value = "涨跌幅 %"
$wsEtfStrong.Range("E1").Value = "RSI"
$wsEtfStrong.Range("F1").Value = "成交量/均量比"
$wsEtfStrong.Range("G1").Value = "策略评分"
$wsEtfStrong.Range("H1").Value = "评级"

# ---------------------------------------------------------------------------
# 4) New sheet "ETF ✅Buy" (headers + SPY / QQQ rows)
# ---------------------------------------------------------------------------
$wsEtfBuy = $wb.Worksheets.Add($null, $wsEtfStrong)
$wsEtfBuy.Name = "ETF ✅Buy"
$wsStockStrong.Range("A1").Copy() | Out-Null
$wsEtfBuy.Range("A1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$wsEtfBuy.Range("A1").Value = "类别"
$wsEtfBuy.Range("B1").Value = "代码"
$wsEtfBuy.Range("C1").Value = "收盘价"
$wsEtfBuy.Range("D1").Value = "涨跌幅 %"
$wsEtfBuy.Range("E1").Value = "RSI"
$wsEtfBuy.Range("F1").Value = "成交量/均量比"
$wsEtfBuy.Range("G1").Value = "策略评分"
$wsEtfBuy.Range("H1").Value = "评级"

$wsEtfBuy.Range("A2").Value = "ETF"
$wsEtfBuy.Range("B2").Value = "SPY"
$wsEtfBuy.Range("C2").Value = 681.4400000000001
$wsEtfBuy.Range("D2").Value = 1.56
$wsEtfBuy.Range("E2").Value = 57.7
$wsEtfBuy.Range("F2").Value = 0.97
$wsEtfBuy.Range("G2").Value = 70
$wsEtfBuy.Range("H2").Value = "✅ 买入"

$wsEtfBuy.Range("A3").Value = "ETF"
$wsEtfBuy.Range("B3").Value = "QQQ"
$wsEtfBuy.Range("C3").Value = 623.23
$wsEtfBuy.Range("D3").Value = 2.21
$wsEtfBuy.Range("E3").Value = 56
$wsEtfBuy.Range("F3").Value = 0.9399999999999999
$wsEtfBuy.Range("G3").Value = 70
$wsEtfBuy.Range("H3").Value = "✅ 买入"

# ---------------------------------------------------------------------------
# 5) "Industry Summary": move after "ETF ✅Buy"; refresh the "✅ 买入" stats
#    (count 5 -> 7, mean 74 -> 72.85714285714286) now that ETFs are included.
# ---------------------------------------------------------------------------
$wsIndustry = $wb.Worksheets.Item("Industry Summary")
$wsIndustry.Move($null, $wsEtfBuy)
$wsIndustry.Range("B2").Value = 7
$wsIndustry.Range("C2").Value = 72.85714285714286

# ---------------------------------------------------------------------------
# 6) New sheet "ETF Overview" (after "Industry Summary")
# ---------------------------------------------------------------------------
$wsEtfOverview = $wb.Worksheets.Add($null, $wsIndustry)
$wsEtfOverview.Name = "ETF Overview"

$wsStockStrong.Range("A1").Copy() | Out-Null
$wsEtfOverview.Range("A1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$wsEtfOverview.Range("A1").Value = "ETF"
$wsEtfOverview.Range("B1").Value = "收盘价"
$wsEtfOverview.Range("C1").Value = "RSI"
$wsEtfOverview.Range("D1").Value = "站上MA20"
$wsEtfOverview.Range("E1").Value = "站上MA50"
$wsEtfOverview.Range("F1").Value = "MACD>Signal"
$wsEtfOverview.Range("G1").Value = "MA20上升"
$wsEtfOverview.Range("H1").Value = "MA50上升"
$wsEtfOverview.Range("I1").Value = "与MA20偏离%"
$wsEtfOverview.Range("J1").Value = "与MA50偏离%"

$wsEtfOverview.Range("A2").Value = "SPY"
$wsEtfOverview.Range("B2").Value = 681.4400000000001
$wsEtfOverview.Range("C2").Value = 57.7
$wsEtfOverview.Range("D2").Value = $true
$wsEtfOverview.Range("E2").Value = $true
$wsEtfOverview.Range("F2").Value = $false
$wsEtfOverview.Range("G2").Value = $true
$wsEtfOverview.Range("H2").Value = $true
$wsEtfOverview.Range("I2").Value = 1.01
$wsEtfOverview.Range("J2").Value = 2.34

$wsEtfOverview.Range("A3").Value = "QQQ"
$wsEtfOverview.Range("B3").Value = 623.23
$wsEtfOverview.Range("C3").Value = 56
$wsEtfOverview.Range("D3").Value = $true
$wsEtfOverview.Range("E3").Value = $true
$wsEtfOverview.Range("F3").Value = $false
$wsEtfOverview.Range("G3").Value = $true
$wsEtfOverview.Range("H3").Value = $true
$wsEtfOverview.Range("I3").Value = 1.08
$wsEtfOverview.Range("J3").Value = 3.46

$wsEtfOverview.Range("A4").Value = "IWM"
$wsEtfOverview.Range("B4").Value = 244.03
$wsEtfOverview.Range("C4").Value = 49.9
$wsEtfOverview.Range("D4").Value = $false
$wsEtfOverview.Range("E4").Value = $true
$wsEtfOverview.Range("F4").Value = $false
$wsEtfOverview.Range("G4").Value = $false
$wsEtfOverview.Range("H4").Value = $true
$wsEtfOverview.Range("I4").Value = -0.71
$wsEtfOverview.Range("J4").Value = 0.47
